$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update conversion text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$oldText = $wsHoja1.Range("A1").Value()
$newText = $oldText.Replace("✅ 1000 Bs = 6.53 = 25905.65 pesos", "✅ 1000 Bs = 6.56 = 25973.77 pesos")
$newText = $newText.Replace("✅ 25905.65 pesos = 6.51 = 970.37 Bs", "✅ 25973.77 pesos = 6.53 = 960.54 Bs")
$wsHoja1.Range("A1").Value = $newText

# --- Sheet "tasas": update rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 152.5
$wsTasas.Range("O10").Value = 3961
$wsTasas.Range("N12").Value = 3975
$wsTasas.Range("O12").Value = 147
